$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Fill in missing "date started" / "date finished" / "number included
# in dataset" values that were filled in during this review pass.
# ------------------------------------------------------------------
# Row 2
$ws.Range("D2").NumberFormat = "d-mmm"
$ws.Range("D2").Value = 43417
$ws.Range("G2").Value = 4

# Row 4
$ws.Range("D4").NumberFormat = "d-mmm"
$ws.Range("D4").Value = 43417
$ws.Range("G4").Value = 1

# Row 5
$ws.Range("D5").NumberFormat = "d-mmm"
$ws.Range("D5").Value = 43417
$ws.Range("G5").Value = 4

# Row 7
$ws.Range("D7").NumberFormat = "d-mmm"
$ws.Range("D7").Value = 43417
$ws.Range("G7").Value = 1

# Row 10
$ws.Range("D10").NumberFormat = "d-mmm"
$ws.Range("D10").Value = 43417
$ws.Range("G10").Value = 1

# Row 13
$ws.Range("D13").NumberFormat = "d-mmm"
$ws.Range("D13").Value = 43417
$ws.Range("G13").Value = 1

# Row 26
$ws.Range("D26").NumberFormat = "d-mmm"
$ws.Range("D26").Value = 43417
$ws.Range("G26").Value = 2

# Row 37
$ws.Range("D37").NumberFormat = "d-mmm"
$ws.Range("D37").Value = 43417
$ws.Range("G37").Value = 1

# Row 40
$ws.Range("D40").NumberFormat = "d-mmm"
$ws.Range("D40").Value = 43417
$ws.Range("G40").Value = 1

# Row 42
$ws.Range("D42").NumberFormat = "d-mmm"
$ws.Range("D42").Value = 43417
$ws.Range("G42").Value = 2

# Row 43
$ws.Range("D43").NumberFormat = "d-mmm"
$ws.Range("D43").Value = 43417
$ws.Range("G43").Value = 1

# Row 54
$ws.Range("D54").NumberFormat = "d-mmm"
$ws.Range("D54").Value = 43417
$ws.Range("E54").NumberFormat = "d-mmm"
$ws.Range("E54").Value = 43417
$ws.Range("G54").Value = 0

# Row 55
$ws.Range("D55").NumberFormat = "d-mmm"
$ws.Range("D55").Value = 43417
$ws.Range("E55").NumberFormat = "d-mmm"
$ws.Range("E55").Value = 43417
$ws.Range("G55").Value = 0

# Row 56
$ws.Range("D56").NumberFormat = "d-mmm"
$ws.Range("D56").Value = 43417
$ws.Range("E56").NumberFormat = "d-mmm"
$ws.Range("E56").Value = 43417
$ws.Range("G56").Value = 0

# Row 57
$ws.Range("D57").NumberFormat = "d-mmm"
$ws.Range("D57").Value = 43417
$ws.Range("E57").NumberFormat = "d-mmm"
$ws.Range("E57").Value = 43417
$ws.Range("G57").Value = 0

# Row 58
$ws.Range("D58").NumberFormat = "d-mmm"
$ws.Range("D58").Value = 43417
$ws.Range("E58").NumberFormat = "d-mmm"
$ws.Range("E58").Value = 43417
$ws.Range("G58").Value = 0

# Row 59
$ws.Range("D59").NumberFormat = "d-mmm"
$ws.Range("D59").Value = 43417
$ws.Range("E59").NumberFormat = "d-mmm"
$ws.Range("E59").Value = 43417
$ws.Range("G59").Value = 1

# Row 60
$ws.Range("D60").NumberFormat = "d-mmm"
$ws.Range("D60").Value = 43417
$ws.Range("E60").NumberFormat = "d-mmm"
$ws.Range("E60").Value = 43417
$ws.Range("G60").Value = 0

# Row 61
$ws.Range("D61").NumberFormat = "d-mmm"
$ws.Range("D61").Value = 43417
$ws.Range("E61").NumberFormat = "d-mmm"
$ws.Range("E61").Value = 43417
$ws.Range("G61").Value = 1

# Row 62
$ws.Range("D62").NumberFormat = "d-mmm"
$ws.Range("D62").Value = 43417
$ws.Range("E62").NumberFormat = "d-mmm"
$ws.Range("E62").Value = 43417
$ws.Range("G62").Value = 0

# Row 63
$ws.Range("D63").NumberFormat = "d-mmm"
$ws.Range("D63").Value = 43417
$ws.Range("E63").NumberFormat = "d-mmm"
$ws.Range("E63").Value = 43417
$ws.Range("G63").Value = 0

# Row 64
$ws.Range("D64").NumberFormat = "d-mmm"
$ws.Range("D64").Value = 43417
$ws.Range("E64").NumberFormat = "d-mmm"
$ws.Range("E64").Value = 43417
$ws.Range("G64").Value = 0

# Row 65
$ws.Range("D65").NumberFormat = "d-mmm"
$ws.Range("D65").Value = 43417
$ws.Range("E65").NumberFormat = "d-mmm"
$ws.Range("E65").Value = 43417
$ws.Range("G65").Value = 0

# Row 66
$ws.Range("D66").NumberFormat = "d-mmm"
$ws.Range("D66").Value = 43417
$ws.Range("E66").NumberFormat = "d-mmm"
$ws.Range("E66").Value = 43417
$ws.Range("G66").Value = 0

# Row 67
$ws.Range("D67").NumberFormat = "d-mmm"
$ws.Range("D67").Value = 43417
$ws.Range("E67").NumberFormat = "d-mmm"
$ws.Range("E67").Value = 43417
$ws.Range("G67").Value = 0

# Row 68
$ws.Range("D68").NumberFormat = "d-mmm"
$ws.Range("D68").Value = 43417
$ws.Range("E68").NumberFormat = "d-mmm"
$ws.Range("E68").Value = 43417
$ws.Range("G68").Value = 0

# Row 69
$ws.Range("D69").NumberFormat = "d-mmm"
$ws.Range("D69").Value = 43417
$ws.Range("E69").NumberFormat = "d-mmm"
$ws.Range("E69").Value = 43417
$ws.Range("G69").Value = 0

# Row 70
$ws.Range("D70").NumberFormat = "d-mmm"
$ws.Range("D70").Value = 43417
$ws.Range("E70").NumberFormat = "d-mmm"
$ws.Range("E70").Value = 43417
$ws.Range("G70").Value = 0

# Row 71
$ws.Range("D71").NumberFormat = "d-mmm"
$ws.Range("D71").Value = 43417
$ws.Range("E71").NumberFormat = "d-mmm"
$ws.Range("E71").Value = 43417
$ws.Range("G71").Value = 1

# ------------------------------------------------------------------
# Restore the scroll position / active selection recorded in the
# workbook when it was last saved.
# ------------------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 31
$ws.Range("D53").Select()
